$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("eil51.tsp")
$ws2 = $wb.Worksheets.Item("berlin52.tsp")

# --- Update existing rows (eil51.tsp) with new experiment values ---
$ws1.Cells.Item(2, 3).Value = 0.7527199268341065
$ws1.Cells.Item(2, 4).Value = 1690.69556594935
$ws1.Cells.Item(2, 5).Value = 587.2384758832502
$ws1.Cells.Item(2, 6).Value = 65.26645673471631
$ws1.Cells.Item(3, 3).Value = 0.7605681657791138
$ws1.Cells.Item(3, 4).Value = 1606.507632883173
$ws1.Cells.Item(3, 5).Value = 601.6161825762687
$ws1.Cells.Item(3, 6).Value = 62.55130257323721
$ws1.Cells.Item(4, 3).Value = 0.7604927301406861
$ws1.Cells.Item(4, 4).Value = 1656.525734252777
$ws1.Cells.Item(4, 5).Value = 582.0750347067585
$ws1.Cells.Item(4, 6).Value = 64.86169682300044
$ws1.Cells.Item(5, 3).Value = 0.7391364812850952
$ws1.Cells.Item(5, 4).Value = 1669.127478383141
$ws1.Cells.Item(5, 5).Value = 629.1458541549424
$ws1.Cells.Item(5, 6).Value = 62.3069021208382
$ws1.Cells.Item(6, 3).Value = 0.7438838243484497
$ws1.Cells.Item(6, 4).Value = 1594.671085733298
$ws1.Cells.Item(6, 5).Value = 633.4624096481157
$ws1.Cells.Item(6, 6).Value = 60.27629676643803
$ws1.Cells.Item(7, 3).Value = 0.7075989246368408
$ws1.Cells.Item(7, 4).Value = 1694.54631786108
$ws1.Cells.Item(7, 5).Value = 641.2655363197239
$ws1.Cells.Item(7, 6).Value = 62.15709599905459
$ws1.Cells.Item(8, 3).Value = 0.009421491622924804
$ws1.Cells.Item(8, 4).Value = 1647.523712337971
$ws1.Cells.Item(8, 5).Value = 1053.532215927919
$ws1.Cells.Item(8, 6).Value = 36.05359315691607
$ws1.Cells.Item(9, 3).Value = 0.006155323982238769
$ws1.Cells.Item(9, 4).Value = 1707.036973975203
$ws1.Cells.Item(9, 5).Value = 1059.096032680611
$ws1.Cells.Item(9, 6).Value = 37.95705372366496
$ws1.Cells.Item(10, 3).Value = 0.005618906021118164
$ws1.Cells.Item(10, 4).Value = 1632.820212375709
$ws1.Cells.Item(10, 5).Value = 1146.085119127346
$ws1.Cells.Item(10, 6).Value = 29.80947256527263
$ws1.Cells.Item(11, 3).Value = 7.845069527626038
$ws1.Cells.Item(11, 4).Value = 1599.786781328893
$ws1.Cells.Item(11, 5).Value = 564.2070818171729
$ws1.Cells.Item(11, 6).Value = 64.73235756151806
$ws1.Cells.Item(12, 3).Value = 7.759456515312195
$ws1.Cells.Item(12, 4).Value = 1680.909565873136
$ws1.Cells.Item(12, 5).Value = 578.7020480937815
$ws1.Cells.Item(12, 6).Value = 65.57208907350235
$ws1.Cells.Item(13, 3).Value = 7.626270771026611
$ws1.Cells.Item(13, 4).Value = 1695.828535289568
$ws1.Cells.Item(13, 5).Value = 570.2023307639324
$ws1.Cells.Item(13, 6).Value = 66.376180203468
$ws1.Cells.Item(14, 1).Value = "optDistCircularIC"
$ws1.Cells.Item(14, 2).Value = 100
$ws1.Cells.Item(14, 3).Value = 0.002634811401367188
$ws1.Cells.Item(14, 4).Value = 1697.38957176831
$ws1.Cells.Item(14, 5).Value = 1197.669716218704
$ws1.Cells.Item(14, 6).Value = 29.44049285215102
$ws1.Cells.Item(15, 1).Value = "optDistCircularIC"
$ws1.Cells.Item(15, 2).Value = 1000
$ws1.Cells.Item(15, 3).Value = 0.02549126148223877
$ws1.Cells.Item(15, 4).Value = 1671.500762575172
$ws1.Cells.Item(15, 5).Value = 761.8136826246156
$ws1.Cells.Item(15, 6).Value = 54.4233721167474
$ws1.Cells.Item(16, 1).Value = "optDistCircularIC"
$ws1.Cells.Item(16, 2).Value = 10000
$ws1.Cells.Item(16, 3).Value = 0.2539310455322266
$ws1.Cells.Item(16, 4).Value = 1640.326677750452
$ws1.Cells.Item(16, 5).Value = 569.2512753969586
$ws1.Cells.Item(16, 6).Value = 65.29646910470106

# --- Update existing rows (berlin52.tsp) with new experiment values ---
$ws2.Cells.Item(2, 3).Value = 0.8538057327270507
$ws2.Cells.Item(2, 4).Value = 30152.46355144759
$ws2.Cells.Item(2, 5).Value = 10562.62370460025
$ws2.Cells.Item(2, 6).Value = 64.96928456085256
$ws2.Cells.Item(3, 3).Value = 0.9253467798233033
$ws2.Cells.Item(3, 4).Value = 30137.48415372318
$ws2.Cells.Item(3, 5).Value = 10476.30416239182
$ws2.Cells.Item(3, 6).Value = 65.23829225772457
$ws2.Cells.Item(4, 3).Value = 0.8165324926376343
$ws2.Cells.Item(4, 4).Value = 30539.66940284299
$ws2.Cells.Item(4, 5).Value = 10398.78358613187
$ws2.Cells.Item(4, 6).Value = 65.94991435904728
$ws2.Cells.Item(5, 3).Value = 0.9262724637985229
$ws2.Cells.Item(5, 4).Value = 30507.73393482665
$ws2.Cells.Item(5, 5).Value = 10805.075173917
$ws2.Cells.Item(5, 6).Value = 64.58250489203895
$ws2.Cells.Item(6, 3).Value = 0.911911416053772
$ws2.Cells.Item(6, 4).Value = 30277.26512266933
$ws2.Cells.Item(6, 5).Value = 11100.15033344575
$ws2.Cells.Item(6, 6).Value = 63.3383322817529
$ws2.Cells.Item(7, 3).Value = 0.90787672996521
$ws2.Cells.Item(7, 4).Value = 29982.78027562029
$ws2.Cells.Item(7, 5).Value = 10917.22045157333
$ws2.Cells.Item(7, 6).Value = 63.58836521758331
$ws2.Cells.Item(8, 3).Value = 0.005371356010437011
$ws2.Cells.Item(8, 4).Value = 29682.95006901009
$ws2.Cells.Item(8, 5).Value = 20363.96588755763
$ws2.Cells.Item(8, 6).Value = 31.39507414117091
$ws2.Cells.Item(9, 3).Value = 0.004100418090820313
$ws2.Cells.Item(9, 4).Value = 29632.53728529375
$ws2.Cells.Item(9, 5).Value = 21742.72220438267
$ws2.Cells.Item(9, 6).Value = 26.62551304652099
$ws2.Cells.Item(10, 3).Value = 0.008357238769531251
$ws2.Cells.Item(10, 4).Value = 29863.53053316114
$ws2.Cells.Item(10, 5).Value = 18516.14271628337
$ws2.Cells.Item(10, 6).Value = 37.99747589883042
$ws2.Cells.Item(11, 3).Value = 9.090412235260009
$ws2.Cells.Item(11, 4).Value = 29580.71606880549
$ws2.Cells.Item(11, 5).Value = 10127.51129329761
$ws2.Cells.Item(11, 6).Value = 65.76313004140684
$ws2.Cells.Item(12, 3).Value = 9.186949181556702
$ws2.Cells.Item(12, 4).Value = 30192.33527117213
$ws2.Cells.Item(12, 5).Value = 10022.06245095608
$ws2.Cells.Item(12, 6).Value = 66.80593812653761
$ws2.Cells.Item(13, 3).Value = 9.017185378074647
$ws2.Cells.Item(13, 4).Value = 29719.36782169909
$ws2.Cells.Item(13, 5).Value = 10105.0933784091
$ws2.Cells.Item(13, 6).Value = 65.99828960348529
$ws2.Cells.Item(14, 1).Value = "optDistCircularIC"
$ws2.Cells.Item(14, 2).Value = 100
$ws2.Cells.Item(14, 3).Value = 0.00268561840057373
$ws2.Cells.Item(14, 4).Value = 29781.35333994811
$ws2.Cells.Item(14, 5).Value = 20982.49073317652
$ws2.Cells.Item(14, 6).Value = 29.54487160584799
$ws2.Cells.Item(15, 1).Value = "optDistCircularIC"
$ws2.Cells.Item(15, 2).Value = 1000
$ws2.Cells.Item(15, 3).Value = 0.02612676620483399
$ws2.Cells.Item(15, 4).Value = 30520.52860017783
$ws2.Cells.Item(15, 5).Value = 13595.46777012297
$ws2.Cells.Item(15, 6).Value = 55.45467790474719
$ws2.Cells.Item(16, 1).Value = "optDistCircularIC"
$ws2.Cells.Item(16, 2).Value = 10000
$ws2.Cells.Item(16, 3).Value = 0.2601830959320068
$ws2.Cells.Item(16, 4).Value = 28944.99884145063
$ws2.Cells.Item(16, 5).Value = 10457.27042856887
$ws2.Cells.Item(16, 6).Value = 63.87192659481624
